$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Upload File")
[void]$ws.Columns("K").Select()
$ws.Columns("K").Delete()
